$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.417.85"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.564.50"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.787.04"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.566.79"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "27.418.41"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "1.361.74"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "1.699.69"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "0.0₇0993"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0493"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
